$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# 1. Update row 9 (history for a player) - new credit history entries
#    and updated aggregate counters.
# -------------------------------------------------------------------

# Remove the I9 cell entirely (stat history column no longer populated
# for this row).
$ws.Range("I9").Clear() | Out-Null

# Extend the history strings held in F9/G9/H9 with the new entries.
$ws.Range("F9").Value = ";0;0;1;1;0;1;0;0;0;0;0;0;0;0;0;0;1;1;0;0;0;0"
$ws.Range("G9").Value = ";14;33;13;13;10;13;31;0;4;42;0;0;0;0;0;3.0;3.0;43;31;3;44;0"
$ws.Range("H9").Value = ";-120.0;-120;+50.0;+100;-100;+1499700;-1499700;-1499700;-100;-100;-100;-120;-100;-100;-100;-4999999500.0;+4999999500.0;+4999999100.0;-14999997000;-120;-120;-120"

# Updated aggregate counters for this user row.
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 17

# -------------------------------------------------------------------
# 2. Header row: center-align all header cells (button/UI polish).
# -------------------------------------------------------------------
$ws.Range("A1:I1").HorizontalAlignment = -4108

# -------------------------------------------------------------------
# 3. Column width / layout adjustments.
# -------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 12.5
$ws.Columns("G").ColumnWidth = 32.5
$ws.Columns("H").ColumnWidth = 75.66666666666667
$ws.Columns("I").ColumnWidth = 22.333333333333336

# -------------------------------------------------------------------
# 4. Selection moves to B9.
# -------------------------------------------------------------------
$ws.Range("B9").Select() | Out-Null
